$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 908253.3
$ws.Range("I17").Value = 409.26923
$ws.Range("J17").Value = 1345363.5
$ws.Range("K17").Value = 1227.80769
$ws.Range("L17").Value = 4036090.5
$ws.Range("M17").Value = -1059.80769
$ws.Range("N17").Value = -4036426.5
# Row 28
$ws.Range("H28").Value = 9014.632
$ws.Range("I28").Value = 1138.3636
$ws.Range("J28").Value = 19844.5
$ws.Range("K28").Value = 1138.3636
$ws.Range("L28").Value = 19844.5
$ws.Range("M28").Value = -653.3635999999999
$ws.Range("N28").Value = -20814.5
# Row 112
$ws.Range("H112").Value = 11479.583
$ws.Range("J112").Value = 11479.583
$ws.Range("L112").Value = 34438.749
$ws.Range("N112").Value = -36654.749
# Row 132
$ws.Range("H132").Value = 1728.2208
$ws.Range("I132").Value = 1549.6428
$ws.Range("K132").Value = 4648.928400000001
$ws.Range("M132").Value = -2118.928400000001
# Row 137
$ws.Range("H137").Value = 1142.0513
$ws.Range("I137").Value = 984.7292
$ws.Range("J137").Value = 1393.7667
$ws.Range("K137").Value = 2954.1876
$ws.Range("L137").Value = 4181.300099999999
$ws.Range("M137").Value = -404.1876000000002
$ws.Range("N137").Value = -9281.3001
# Row 138
$ws.Range("H138").Value = 1068.71
$ws.Range("I138").Value = 556.96075
$ws.Range("J138").Value = 1601.3469
$ws.Range("K138").Value = 1670.88225
$ws.Range("L138").Value = 4804.0407
$ws.Range("M138").Value = 3469.11775
$ws.Range("N138").Value = -15084.0407
# Row 141
$ws.Range("H141").Value = 2522.4363
$ws.Range("I141").Value = 874.1111
$ws.Range("J141").Value = 9939.9
$ws.Range("K141").Value = 2622.3333
$ws.Range("L141").Value = 29819.7
$ws.Range("M141").Value = 2557.6667
$ws.Range("N141").Value = -40179.7

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 789163.3
$ws.Range("I32").Value = 939281.6
$ws.Range("J32").Value = 15476.615
$ws.Range("K32").Value = 939281.6
$ws.Range("L32").Value = 15476.615
$ws.Range("M32").Value = -938994.6
$ws.Range("N32").Value = -16050.615
# Row 45
$ws.Range("H45").Value = 3464.889
$ws.Range("I45").Value = 3499.1428
$ws.Range("J45").Value = 3345
$ws.Range("K45").Value = 3499.1428
$ws.Range("L45").Value = 3345
$ws.Range("M45").Value = -3122.1428
$ws.Range("N45").Value = -4099
# Row 88
$ws.Range("H88").Value = 1555.5555
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 2000
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -2812
# Row 91
$ws.Range("H91").Value = 1555.5555
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 2000
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 2000
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -4808
# Row 122
$ws.Range("H122").Value = 168674.5
$ws.Range("I122").Value = 201626.6
$ws.Range("K122").Value = 604879.8
$ws.Range("M122").Value = -602429.8
# Row 132
$ws.Range("H132").Value = 3093.913
$ws.Range("I132").Value = 3167.4443
$ws.Range("K132").Value = 9502.332900000001
$ws.Range("M132").Value = -6972.332900000001

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 15626162
$ws.Range("I105").Value = 17858272
$ws.Range("J105").Value = 1400
$ws.Range("K105").Value = 17858272
$ws.Range("L105").Value = 1400
$ws.Range("M105").Value = -17856525
$ws.Range("N105").Value = -4894
# Row 134
$ws.Range("H134").Value = 3425.5356
$ws.Range("I134").Value = 3200.75
$ws.Range("J134").Value = 3987.5
$ws.Range("K134").Value = 9602.25
$ws.Range("L134").Value = 11962.5
$ws.Range("M134").Value = -7067.25
$ws.Range("N134").Value = -17032.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4608.067
$ws.Range("I31").Value = 1413.6207
$ws.Range("K31").Value = 1413.6207
$ws.Range("M31").Value = -1118.6207
# Row 34
$ws.Range("H34").Value = 4608.067
$ws.Range("I34").Value = 1413.6207
$ws.Range("K34").Value = 1413.6207
$ws.Range("M34").Value = -1211.6207
# Row 132
$ws.Range("H132").Value = 4903817.5
$ws.Range("I132").Value = 1666.238
$ws.Range("J132").Value = 12822677
$ws.Range("K132").Value = 4998.714
$ws.Range("L132").Value = 38468031
$ws.Range("M132").Value = -2468.714
$ws.Range("N132").Value = -38473091
# Row 141
$ws.Range("H141").Value = 76798.64
$ws.Range("J141").Value = 71665.25
$ws.Range("L141").Value = 71665.25
$ws.Range("N141").Value = -82025.25

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 3140.698
$ws.Range("I131").Value = 395.86667
$ws.Range("J131").Value = 4224.184
$ws.Range("K131").Value = 1187.60001
$ws.Range("L131").Value = 12672.552
$ws.Range("M131").Value = 3852.39999
$ws.Range("N131").Value = -22752.552

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 9304.083000000001
$ws.Range("I70").Value = 10142
$ws.Range("K70").Value = 10142
$ws.Range("M70").Value = -9872
# Row 73
$ws.Range("H73").Value = 9304.083000000001
$ws.Range("I73").Value = 10142
$ws.Range("K73").Value = 10142
$ws.Range("M73").Value = -9206
# Row 132
$ws.Range("H132").Value = 2670.0652
$ws.Range("I132").Value = 2119.484
$ws.Range("J132").Value = 3807.9333
$ws.Range("K132").Value = 6358.451999999999
$ws.Range("L132").Value = 11423.7999
$ws.Range("M132").Value = -3828.451999999999
$ws.Range("N132").Value = -16483.7999

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 927
$ws.Range("I46").Value = 629.1667
$ws.Range("J46").Value = 1284.4
$ws.Range("K46").Value = 629.1667
$ws.Range("L46").Value = 1284.4
$ws.Range("M46").Value = -441.1667
$ws.Range("N46").Value = -1660.4
# Row 122
$ws.Range("H122").Value = 4373.654
$ws.Range("I122").Value = 3368.5715
$ws.Range("J122").Value = 4743.9473
$ws.Range("K122").Value = 10105.7145
$ws.Range("L122").Value = 14231.8419
$ws.Range("M122").Value = -7655.7145
$ws.Range("N122").Value = -19131.8419
# Row 132
$ws.Range("H132").Value = 2156.6365
$ws.Range("I132").Value = 1957.7333
$ws.Range("J132").Value = 2858.647
$ws.Range("K132").Value = 5873.199900000001
$ws.Range("L132").Value = 8575.940999999999
$ws.Range("M132").Value = -3343.199900000001
$ws.Range("N132").Value = -13635.941
# Row 136
$ws.Range("H136").Value = 5052314
$ws.Range("I136").Value = 1789
$ws.Range("K136").Value = 5367
$ws.Range("M136").Value = -2817

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2881.125
$ws.Range("I122").Value = 2701.4814
$ws.Range("J122").Value = 3254.2307
$ws.Range("K122").Value = 8104.4442
$ws.Range("L122").Value = 9762.6921
$ws.Range("M122").Value = -5654.4442
$ws.Range("N122").Value = -14662.6921
# Row 126
$ws.Range("H126").Value = 1194.8334
$ws.Range("I126").Value = 1235.5
$ws.Range("J126").Value = 1052.5
$ws.Range("K126").Value = 3706.5
$ws.Range("L126").Value = 3157.5
$ws.Range("M126").Value = -1236.5
$ws.Range("N126").Value = -8097.5
# Row 136
$ws.Range("H136").Value = 2092.8452
$ws.Range("I136").Value = 2119.8867
$ws.Range("J136").Value = 2046.6129
$ws.Range("K136").Value = 6359.6601
$ws.Range("L136").Value = 6139.8387
$ws.Range("M136").Value = -3809.6601
$ws.Range("N136").Value = -11239.8387
